$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3890.4443
$ws.Range("J19").Value = 2163.8
$ws.Range("L19").Value = 2163.8
$ws.Range("N19").Value = -2513.8
$ws.Range("H32").Value = 12964.833
$ws.Range("I32").Value = 12521.667
$ws.Range("J32").Value = 13186.417
$ws.Range("K32").Value = 12521.667
$ws.Range("L32").Value = 13186.417
$ws.Range("M32").Value = -12195.667
$ws.Range("N32").Value = -13838.417
$ws.Range("H40").Value = 2993.625
$ws.Range("I40").Value = 1790
$ws.Range("K40").Value = 1790
$ws.Range("M40").Value = -1615
$ws.Range("H45").Value = 3400
$ws.Range("I45").Value = 3400
$ws.Range("K45").Value = 10200
$ws.Range("M45").Value = -10008
$ws.Range("H132").Value = 2519778.5
$ws.Range("J132").Value = 1199.6666
$ws.Range("L132").Value = 3598.9998
$ws.Range("N132").Value = -8658.9998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1700
$ws.Range("I41").Value = 1700
$ws.Range("K41").Value = 1700
$ws.Range("M41").Value = -1286
$ws.Range("H45").Value = 2196.3845
$ws.Range("I45").Value = 1090.0769
$ws.Range("J45").Value = 3302.6924
$ws.Range("K45").Value = 1090.0769
$ws.Range("L45").Value = 3302.6924
$ws.Range("M45").Value = -713.0769
$ws.Range("N45").Value = -4056.6924
$ws.Range("H97").Value = 793.5
$ws.Range("I97").Value = 665.7045000000001
$ws.Range("K97").Value = 665.7045000000001
$ws.Range("M97").Value = -169.7045000000001
$ws.Range("H110").Value = 6146.0454
$ws.Range("I110").Value = 6034.6313
$ws.Range("J110").Value = 6851.6665
$ws.Range("K110").Value = 6034.6313
$ws.Range("L110").Value = 6851.6665
$ws.Range("M110").Value = -3989.6313
$ws.Range("N110").Value = -10941.6665
$ws.Range("H122").Value = 1455.9656
$ws.Range("I122").Value = 1274.5834
$ws.Range("J122").Value = 2326.6
$ws.Range("K122").Value = 3823.7502
$ws.Range("L122").Value = 6979.799999999999
$ws.Range("M122").Value = -1373.7502
$ws.Range("N122").Value = -11879.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1499.8572
$ws.Range("I99").Value = 1509.8
$ws.Range("K99").Value = 1509.8
$ws.Range("M99").Value = -11.79999999999995
$ws.Range("H105").Value = 2215.2222
$ws.Range("I105").Value = 2004
$ws.Range("K105").Value = 2004
$ws.Range("M105").Value = -257
$ws.Range("H107").Value = 2289.6775
$ws.Range("I107").Value = 2289.6775
$ws.Range("K107").Value = 2289.6775
$ws.Range("M107").Value = -369.6774999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 670.55554
$ws.Range("J22").Value = 753.4286
$ws.Range("L22").Value = 753.4286
$ws.Range("N22").Value = -1453.4286
$ws.Range("H31").Value = 3449675.5
$ws.Range("I31").Value = 5556488
$ws.Range("J31").Value = 2163.9092
$ws.Range("K31").Value = 5556488
$ws.Range("L31").Value = 2163.9092
$ws.Range("M31").Value = -5556193
$ws.Range("N31").Value = -2753.9092
$ws.Range("H34").Value = 3449675.5
$ws.Range("I34").Value = 5556488
$ws.Range("J34").Value = 2163.9092
$ws.Range("K34").Value = 5556488
$ws.Range("L34").Value = 2163.9092
$ws.Range("M34").Value = -5556286
$ws.Range("N34").Value = -2567.9092
$ws.Range("H74").Value = 90000
$ws.Range("J74").Value = 90000
$ws.Range("L74").Value = 90000
$ws.Range("N74").Value = -91748
$ws.Range("H77").Value = 90000
$ws.Range("J77").Value = 90000
$ws.Range("L77").Value = 270000
$ws.Range("N77").Value = -278736
$ws.Range("H94").Value = 1896.5264
$ws.Range("I94").Value = 1339.1
$ws.Range("J94").Value = 2515.889
$ws.Range("K94").Value = 1339.1
$ws.Range("L94").Value = 2515.889
$ws.Range("M94").Value = -888.0999999999999
$ws.Range("N94").Value = -3417.889
$ws.Range("H105").Value = 862.5833
$ws.Range("I105").Value = 807.2857
$ws.Range("K105").Value = 807.2857
$ws.Range("M105").Value = 939.7143
$ws.Range("H107").Value = 823.65
$ws.Range("I107").Value = 557.2941
$ws.Range("K107").Value = 557.2941
$ws.Range("M107").Value = 1362.7059
$ws.Range("H134").Value = 2510.3235
$ws.Range("I134").Value = 2229.3447
$ws.Range("J134").Value = 4140
$ws.Range("K134").Value = 6688.034100000001
$ws.Range("L134").Value = 12420
$ws.Range("M134").Value = -4153.034100000001
$ws.Range("N134").Value = -17490

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4526.1113
$ws.Range("I69").Value = 985
$ws.Range("K69").Value = 2955
$ws.Range("M69").Value = -2144
$ws.Range("H72").Value = 4526.1113
$ws.Range("I72").Value = 985
$ws.Range("K72").Value = 8865
$ws.Range("M72").Value = -4809
$ws.Range("H137").Value = 2593.25
$ws.Range("I137").Value = 2641.4546
$ws.Range("J137").Value = 2534.3333
$ws.Range("K137").Value = 7924.3638
$ws.Range("L137").Value = 7602.999899999999
$ws.Range("M137").Value = -2824.3638
$ws.Range("N137").Value = -17802.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 21422.88
$ws.Range("I102").Value = 24920.096
$ws.Range("J102").Value = 3062.5
$ws.Range("K102").Value = 24920.096
$ws.Range("L102").Value = 3062.5
$ws.Range("M102").Value = -23298.096
$ws.Range("N102").Value = -6306.5
$ws.Range("H122").Value = 3181.9
$ws.Range("I122").Value = 2874.9565
$ws.Range("K122").Value = 8624.869499999999
$ws.Range("M122").Value = -6174.869499999999
$ws.Range("H132").Value = 2266.2856
$ws.Range("I132").Value = 1859.5
$ws.Range("J132").Value = 4707
$ws.Range("K132").Value = 5578.5
$ws.Range("L132").Value = 14121
$ws.Range("M132").Value = -3048.5
$ws.Range("N132").Value = -19181

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4249.5713
$ws.Range("J7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3224
$ws.Range("H46").Value = 3217.5789
$ws.Range("I46").Value = 878.3333
$ws.Range("J46").Value = 4297.231
$ws.Range("K46").Value = 878.3333
$ws.Range("L46").Value = 4297.231
$ws.Range("M46").Value = -690.3333
$ws.Range("N46").Value = -4673.231
$ws.Range("H61").Value = 2523.2222
$ws.Range("I61").Value = 2558.5
$ws.Range("K61").Value = 2558.5
$ws.Range("M61").Value = -2356.5
$ws.Range("H68").Value = 2184.8
$ws.Range("J68").Value = 2422.4614
$ws.Range("L68").Value = 2422.4614
$ws.Range("N68").Value = -3920.4614
$ws.Range("H71").Value = 2184.8
$ws.Range("J71").Value = 2422.4614
$ws.Range("L71").Value = 12112.307
$ws.Range("N71").Value = -19600.307
$ws.Range("H88").Value = 53247.418
$ws.Range("J88").Value = 59996.375
$ws.Range("L88").Value = 59996.375
$ws.Range("N88").Value = -60852.375
$ws.Range("H91").Value = 53247.418
$ws.Range("J91").Value = 59996.375
$ws.Range("L91").Value = 59996.375
$ws.Range("N91").Value = -62960.375
$ws.Range("H93").Value = 1274.2667
$ws.Range("J93").Value = 2555
$ws.Range("L93").Value = 2555
$ws.Range("N93").Value = -5051
$ws.Range("H113").Value = 2523.2222
$ws.Range("I113").Value = 2558.5
$ws.Range("K113").Value = 2558.5
$ws.Range("M113").Value = -388.5
$ws.Range("H122").Value = 3999.2
$ws.Range("I122").Value = 3999
$ws.Range("K122").Value = 11997
$ws.Range("M122").Value = -9547
$ws.Range("H126").Value = 4249.5713
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 3871.4583
$ws.Range("I132").Value = 3926.353
$ws.Range("J132").Value = 3738.1428
$ws.Range("K132").Value = 11779.059
$ws.Range("L132").Value = 11214.4284
$ws.Range("M132").Value = -9249.059000000001
$ws.Range("N132").Value = -16274.4284

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1198.3636
$ws.Range("I107").Value = 1180.25
$ws.Range("K107").Value = 3540.75
$ws.Range("M107").Value = -1620.75
$ws.Range("H122").Value = 9285577
$ws.Range("I122").Value = 10028248
$ws.Range("J122").Value = 2195.5
$ws.Range("K122").Value = 30084744
$ws.Range("L122").Value = 6586.5
$ws.Range("M122").Value = -30082294
$ws.Range("N122").Value = -11486.5
$ws.Range("H132").Value = 10901851
$ws.Range("I132").Value = 13196478
$ws.Range("J132").Value = 2374.75
$ws.Range("K132").Value = 39589434
$ws.Range("L132").Value = 7124.25
$ws.Range("M132").Value = -39586904
$ws.Range("N132").Value = -12184.25
